$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 484, shifting rows 484:510 down to 485:511
$ws.Rows.Item(484).Insert()

# Populate the new row 484 with the new data record
$ws.Range("A484").Value = 10
$ws.Range("B484").Value = "Vega Modelo de Temuco"
$ws.Range("C484").Value = "La Araucanía"
$ws.Range("D484").Value = 44516
$ws.Range("D484").NumberFormat = $ws.Range("D485").NumberFormat
$ws.Range("E484").Value = 9
$ws.Range("F484").Value = 100114001
$ws.Range("G484").Value = "Papa"
$ws.Range("H484").Value = "Patagonia"
$ws.Range("I484").Value = "1a nueva(o)"
$ws.Range("J484").Value = 580
$ws.Range("K484").Value = 12000
$ws.Range("L484").Value = 13000
$ws.Range("M484").Value = 12655
$ws.Range("N484").Value = "$/saco 25 kilos"
$ws.Range("O484").Value = "Provincia de Cautín"
$ws.Range("P484").Value = 506
$ws.Range("Q484").Value = 25
$ws.Range("R484").Value = "Hortaliza"
